# Konami code and counting!
#
# Mark a handful of "Hats" achievements as implemented in column A
# (the "Impl?" column), then leave the workbook's view/selection state
# the way the author left it: scrolled down on the "Hats" sheet with
# A42 selected (and "Hats" the active/tabbed sheet), while "Camos" is
# left with D6 selected.

$wb = $excel.ActiveWorkbook

$wsHats = $wb.Worksheets.Item("Hats")
$wsCamos = $wb.Worksheets.Item("Camos")

# Rows 25 & 26 get marked with the "PIPE" type value (matching column D
# for those two rows); the rest get a plain "x" like the rest of the
# column.
$wsHats.Range("A25").Value = "PIPE"
$wsHats.Range("A26").Value = "PIPE"

$wsHats.Range("A37").Value = "x"
$wsHats.Range("A38").Value = "x"
$wsHats.Range("A39").Value = "x"
$wsHats.Range("A41").Value = "x"
$wsHats.Range("A45").Value = "x"
$wsHats.Range("A46").Value = "x"
$wsHats.Range("A47").Value = "x"

# Camos was the active tab / selection before; move the selection on
# Camos first (without activating it) so it's left on D6...
$wsCamos.Range("D6").Select() | Out-Null

# ...then activate Hats and leave it scrolled/selected at A42.
$wsHats.Activate()
$wsHats.Range("A42").Select() | Out-Null
